{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text,alignment,font/bold,font/size\");\nawait context.sync();\n\n// Locate the title paragraph (\"STAR MORTORS\") and remove the empty,\n// right-aligned / bold / 7pt paragraph that immediately follows it.\nlet titleIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"STAR MORTORS\") {\n    titleIndex = i;\n    break;\n  }\n}\n\nif (titleIndex !== -1 && titleIndex + 1 < paragraphs.items.length) {\n  const candidate = paragraphs.items[titleIndex + 1];\n  if (candidate.text.trim() === \"\") {\n    candidate.delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the title paragraph (\"STAR MORTORS\") and remove the empty,\n# right-aligned / bold / 7pt paragraph that immediately follows it.\n$titlePara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"STAR MORTORS\") {\n        $titlePara = $p\n        break\n    }\n}\n\nif ($titlePara -ne $null) {\n    $nextPara = $titlePara.Next()\n    if ($nextPara -ne $null) {\n        $nextText = $nextPara.Range.Text.TrimEnd([char]13, [char]7)\n        if ($nextText -eq \"\") {\n            $nextPara.Range.Delete()\n        }\n    }\n}\n"}
